$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 1020, shifting the existing rows (old 1020-1072)
# down to 1022-1074. The inserted rows inherit the formatting (incl. the
# date style on column D) from the row they are inserted above, matching
# the target workbook exactly.
$ws.Range("A1020:A1021").EntireRow.Insert()

# New row 1020: "Primera" quality entry for the new weekly date 45041 (2023-04-25)
$ws.Range("A1020").Value = 8
$ws.Range("B1020").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1020").Value = "Coquimbo"
$ws.Range("D1020").Value = 45041
$ws.Range("E1020").Value = 4
$ws.Range("F1020").Value = 100112008
$ws.Range("G1020").Value = "Coliflor"
$ws.Range("H1020").Value = "Sin especificar"
$ws.Range("I1020").Value = "Primera"
$ws.Range("J1020").Value = 2100
$ws.Range("K1020").Value = 800
$ws.Range("L1020").Value = 900
$ws.Range("M1020").Value = 850
$ws.Range("N1020").Value = "$/unidad"
$ws.Range("O1020").Value = "Provincia del Elquí"
$ws.Range("P1020").Value = 850
$ws.Range("Q1020").Value = 1
$ws.Range("R1020").Value = "Hortaliza"

# New row 1021: "Segunda" quality entry for the same new weekly date
$ws.Range("A1021").Value = 8
$ws.Range("B1021").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1021").Value = "Coquimbo"
$ws.Range("D1021").Value = 45041
$ws.Range("E1021").Value = 4
$ws.Range("F1021").Value = 100112008
$ws.Range("G1021").Value = "Coliflor"
$ws.Range("H1021").Value = "Sin especificar"
$ws.Range("I1021").Value = "Segunda"
$ws.Range("J1021").Value = 1400
$ws.Range("K1021").Value = 600
$ws.Range("L1021").Value = 700
$ws.Range("M1021").Value = 650
$ws.Range("N1021").Value = "$/unidad"
$ws.Range("O1021").Value = "Provincia del Elquí"
$ws.Range("P1021").Value = 650
$ws.Range("Q1021").Value = 1
$ws.Range("R1021").Value = "Hortaliza"

$ws.Range("D1020:D1021").NumberFormat = "YYYY-MM-DD HH:MM:SS"
